$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.429188333333333
$ws.Range("H2").Value = 7.287565
$ws.Range("I2").Value = 0.1102134218696762
$ws.Range("J2").Value = 0.110840272037245
$ws.Range("M2").Value = 2.223710666666667
$ws.Range("N2").Value = 6.671132
$ws.Range("O2").Value = 0.529381647492601
$ws.Range("P2").Value = 0.5293816474926011
$ws.Range("Q2").Value = 5.401812008175556
$ws.Range("R2").Value = 48.61630807358
$ws.Range("S2").Value = 0.05834496284516624
$ws.Range("T2").Value = 0.05867680581960485

$ws.Range("G3").Value = 2.429188333333333
$ws.Range("H3").Value = 7.287565
$ws.Range("I3").Value = 0.1102134218696762
$ws.Range("J3").Value = 0.110840272037245
$ws.Range("O3").Value = 0.4706183525073989
$ws.Range("P3").Value = 0.470618352507399
$ws.Range("Q3").Value = 4.802191159975555
$ws.Range("R3").Value = 43.21972043978
$ws.Range("S3").Value = 0.05186845902450993
$ws.Range("T3").Value = 0.05216346621764018

$ws.Range("I4").Value = 0.1288109145567312
$ws.Range("J4").Value = 0.1295435398759056
$ws.Range("M4").Value = 2.223710666666667
$ws.Range("N4").Value = 6.671132
$ws.Range("O4").Value = 0.529381647492601
$ws.Range("P4").Value = 0.5293816474926011
$ws.Range("Q4").Value = 6.313317681574223
$ws.Range("R4").Value = 56.81985913416801
$ws.Range("S4").Value = 0.068190134163071
$ws.Range("T4").Value = 0.06857797256153039

$ws.Range("I5").Value = 0.1288109145567312
$ws.Range("J5").Value = 0.1295435398759056
$ws.Range("O5").Value = 0.4706183525073989
$ws.Range("P5").Value = 0.470618352507399
$ws.Range("S5").Value = 0.06062078039366015
$ws.Range("T5").Value = 0.06096556731437525

$ws.Range("G6").Value = 8.195700333333333
$ws.Range("H6").Value = 24.587101
$ws.Range("I6").Value = 0.3718427945500777
$ws.Range("J6").Value = 0.3739576886720351
$ws.Range("M6").Value = 2.223710666666667
$ws.Range("N6").Value = 6.671132
$ws.Range("O6").Value = 0.529381647492601
$ws.Range("P6").Value = 0.5293816474926011
$ws.Range("Q6").Value = 18.22486625203689
$ws.Range("R6").Value = 164.023796268332
$ws.Range("S6").Value = 0.1968467511871729
$ws.Range("T6").Value = 0.1979663373217271

$ws.Range("G7").Value = 8.195700333333333
$ws.Range("H7").Value = 24.587101
$ws.Range("I7").Value = 0.3718427945500777
$ws.Range("J7").Value = 0.3739576886720351
$ws.Range("O7").Value = 0.4706183525073989
$ws.Range("P7").Value = 0.470618352507399
$ws.Range("Q7").Value = 16.20183958175689
$ws.Range("R7").Value = 145.816556235812
$ws.Range("S7").Value = 0.1749960433629048
$ws.Range("T7").Value = 0.175991351350308

$ws.Range("G8").Value = 0.3739505
$ws.Range("H8").Value = 0.747901
$ws.Range("I8").Value = 0.01696631078345497
$ws.Range("J8").Value = 0.0113752056135249
$ws.Range("M8").Value = 2.223710666666667
$ws.Range("N8").Value = 6.671132
$ws.Range("O8").Value = 0.529381647492601
$ws.Range("P8").Value = 0.5293816474926011
$ws.Range("Q8").Value = 0.8315577156553334
$ws.Range("R8").Value = 4.989346293932
$ws.Range("S8").Value = 0.008981653554416874
$ws.Range("T8").Value = 0.006021825088254897

$ws.Range("G9").Value = 0.3739505
$ws.Range("H9").Value = 0.747901
$ws.Range("I9").Value = 0.01696631078345497
$ws.Range("J9").Value = 0.0113752056135249
$ws.Range("O9").Value = 0.4706183525073989
$ws.Range("P9").Value = 0.470618352507399
$ws.Range("Q9").Value = 0.7392517742353334
$ws.Range("R9").Value = 4.435510645412
$ws.Range("S9").Value = 0.007984657229038097
$ws.Range("T9").Value = 0.005353380525270006

$ws.Range("G10").Value = 8.202836333333334
$ws.Range("H10").Value = 24.608509
$ws.Range("I10").Value = 0.3721665582400601
$ws.Range("J10").Value = 0.3742832938012893
$ws.Range("M10").Value = 2.223710666666667
$ws.Range("N10").Value = 6.671132
$ws.Range("O10").Value = 0.529381647492601
$ws.Range("P10").Value = 0.5293816474926011
$ws.Range("Q10").Value = 18.24073465135422
$ws.Range("R10").Value = 164.166611862188
$ws.Range("S10").Value = 0.197018145742774
$ws.Range("T10").Value = 0.1981387067014838

$ws.Range("G11").Value = 8.202836333333334
$ws.Range("H11").Value = 24.608509
$ws.Range("I11").Value = 0.3721665582400601
$ws.Range("J11").Value = 0.3742832938012893
$ws.Range("O11").Value = 0.4706183525073989
$ws.Range("P11").Value = 0.470618352507399
$ws.Range("Q11").Value = 16.21594653083422
$ws.Range("R11").Value = 145.943518777508
$ws.Range("S11").Value = 0.175148412497286
$ws.Range("T11").Value = 0.175991351350308

